$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp footer (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 08:35"

# --- Countries table: Ucrania jumps ahead of Austria/Japon/Banglades ---
# Before:  row35=Austria, row36=Japon, row37=Banglades, row38=Ucrania
# After:   row35=Ucrania(new numbers), row36=Austria, row37=Japon, row38=Banglades
# (row39=Rumania is untouched)

# Capture the current (pre-edit) data for the three rows that shift down.
$austria   = @($ws.Cells.Item(35,2).Value2, $ws.Cells.Item(35,3).Value2, $ws.Cells.Item(35,4).Value2, $ws.Cells.Item(35,5).Value2, $ws.Cells.Item(35,6).Value2, $ws.Cells.Item(35,7).Value2, $ws.Cells.Item(35,8).Value2)
$japon     = @($ws.Cells.Item(36,2).Value2, $ws.Cells.Item(36,3).Value2, $ws.Cells.Item(36,4).Value2, $ws.Cells.Item(36,5).Value2, $ws.Cells.Item(36,6).Value2, $ws.Cells.Item(36,7).Value2, $ws.Cells.Item(36,8).Value2)
$banglades = @($ws.Cells.Item(37,2).Value2, $ws.Cells.Item(37,3).Value2, $ws.Cells.Item(37,4).Value2, $ws.Cells.Item(37,5).Value2, $ws.Cells.Item(37,6).Value2, $ws.Cells.Item(37,7).Value2, $ws.Cells.Item(37,8).Value2)

# Row 35 becomes Ucrania with its freshly updated figures.
$ws.Cells.Item(35,1).Value = "Ucrania"
$ws.Cells.Item(35,2).Value = 16023
$ws.Cells.Item(35,3).Value = 375
$ws.Cells.Item(35,4).Value = 3373
$ws.Cells.Item(35,5).Value = 12225
$ws.Cells.Item(35,6).Value = 211
$ws.Cells.Item(35,7).Value = 17
$ws.Cells.Item(35,8).Value = 425

# Row 36 becomes Austria, carrying forward its previous (unchanged) numbers.
$ws.Cells.Item(36,1).Value = "Austria"
$ws.Cells.Item(36,2).Value = $austria[0]
$ws.Cells.Item(36,3).Value = $austria[1]
$ws.Cells.Item(36,4).Value = $austria[2]
$ws.Cells.Item(36,5).Value = $austria[3]
$ws.Cells.Item(36,6).Value = $austria[4]
$ws.Cells.Item(36,7).Value = $austria[5]
$ws.Cells.Item(36,8).Value = $austria[6]

# Row 37 becomes Japon, carrying forward its previous (unchanged) numbers.
$ws.Cells.Item(37,1).Value = "Japon"
$ws.Cells.Item(37,2).Value = $japon[0]
$ws.Cells.Item(37,3).Value = $japon[1]
$ws.Cells.Item(37,4).Value = $japon[2]
$ws.Cells.Item(37,5).Value = $japon[3]
$ws.Cells.Item(37,6).Value = $japon[4]
$ws.Cells.Item(37,7).Value = $japon[5]
$ws.Cells.Item(37,8).Value = $japon[6]

# Row 38 becomes Banglades, carrying forward its previous (unchanged) numbers.
$ws.Cells.Item(38,1).Value = "Banglades"
$ws.Cells.Item(38,2).Value = $banglades[0]
$ws.Cells.Item(38,3).Value = $banglades[1]
$ws.Cells.Item(38,4).Value = $banglades[2]
$ws.Cells.Item(38,5).Value = $banglades[3]
$ws.Cells.Item(38,6).Value = $banglades[4]
$ws.Cells.Item(38,7).Value = $banglades[5]
$ws.Cells.Item(38,8).Value = $banglades[6]

# --- Provincias/countries case-count refresh for a couple of other rows ---
# Row 51: Chequia
$ws.Cells.Item(51,2).Value = 8177
$ws.Cells.Item(51,3).Value = 1
$ws.Cells.Item(51,4).Value = 4738
$ws.Cells.Item(51,5).Value = 3156
$ws.Cells.Item(51,6).Value = 42
$ws.Cells.Item(51,7).Value = 1
$ws.Cells.Item(51,8).Value = 283

# Row 119: Georgia
$ws.Cells.Item(119,2).Value = 639
$ws.Cells.Item(119,3).Value = 1
$ws.Cells.Item(119,4).Value = 349
$ws.Cells.Item(119,5).Value = 279
